# Task 9 slide deck update:
#  1. Slide 1 (title slide) subtitle shape "Resources and microservices"
#     gains a second, centered paragraph with the name of the person
#     responsible for the K8s cluster: "Smolkin Mikhail" (sz=28, dk1).
#  2. Slide 2 (Product description) collapses the five separate runs that
#     spelled out "Dandamaev Gadji" back into the single run
#     " Daniel; Dandamaev Gadji; Tsaturyan Konstantin; " that precedes the
#     existing "Smolkin Mikhail" entry.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1 - add "Smolkin Mikhail" as its own centered paragraph under
#    "Resources and microservices".
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Resources and microservices") {
        $subtitle = $shp
    }
}

$subtitleRange = $subtitle.TextFrame.TextRange
# New paragraph break + text; the engine copies the paragraph's pPr
# (centred, no bullet) from the paragraph it is appended after.
$subtitleRange.InsertAfter("`rSmolkin Mikhail") | Out-Null

$fullText = $subtitleRange.Text
$nameStart = $fullText.IndexOf("Smolkin Mikhail") + 1

$firstName = $subtitleRange.Characters($nameStart, "Smolkin".Length)
$firstName.Font.Size = 28
$firstName.Font.Color.ObjectThemeColor = 1   # msoThemeColorDark1 -> schemeClr dk1

$lastName = $subtitleRange.Characters($nameStart + "Smolkin".Length, " Mikhail".Length)
$lastName.Font.Size = 28
$lastName.Font.Color.ObjectThemeColor = 1    # msoThemeColorDark1 -> schemeClr dk1

# ---------------------------------------------------------------------
# 2) Slide 2 - merge the " Daniel; " / "Dandamaev" / " " / "Gadji" /
#    "; Tsaturyan Konstantin; " runs back into one run.
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
for ($i = 1; $i -le $slide2.Shapes.Count; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text.IndexOf("Team K8C") -ge 0) {
        $body = $shp
    }
}

$bodyRange = $body.TextFrame.TextRange
$bodyText = $bodyRange.Text

$mergedText = " Daniel; Dandamaev Gadji; Tsaturyan Konstantin; "
$mergeStart = $bodyText.IndexOf(" Daniel; ") + 1
$oldEndMarker = "; Tsaturyan Konstantin; "
$mergeEnd = $bodyText.IndexOf($oldEndMarker) + $oldEndMarker.Length
$mergeLen = $mergeEnd - $mergeStart + 1

$mergeRange = $bodyRange.Characters($mergeStart, $mergeLen)
$mergeRange.Text = $mergedText
